$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "FORMALIZACIÓN DE LA RESERVA ..." heading: merge the two runs
#    that were split around the old _GoBack bookmark into a single
#    run, and remove the bookmark from this location.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$rng1 = $d.Content
$rng1.Find.Execute(
    "FORMALIZACIÓN DE LA RESERVA (en caso de propuesta aprobada con reserva)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "FORMALIZACIÓN DE LA RESERVA (en caso de propuesta aprobada con reserva)",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Fix the bank account number (typo correction): insert a missing
#    "6" digit into ES50-2100-8981-69-020001661 -> ...-0200016661
#    Track the insertion so Word keeps it as its own run (mirrors
#    the 3-run split in the authored document), then accept the
#    revision so no tracked-change markup remains in the output.
# ------------------------------------------------------------------
$d.TrackRevisions = $true

$rng2 = $d.Content
$rng2.Find.Execute(
    "ES50-2100-8981-69-02000166",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insertPoint = $rng2.Duplicate
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("6")

$d.TrackRevisions = $false
$d.AcceptAllRevisions() | Out-Null

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark: it now belongs right after
#    "... de los compradores. " at the end of the COBRO RESERVA
#    paragraph (a zero-length bookmark at the end of the paragraph,
#    right before the paragraph mark).
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute(
    "d.n.i. de los compradores. ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$endPoint = $rng3.Duplicate
$endPoint.Collapse(0)
$targetPos = $endPoint.Start

# Work around a collapsed-range-at-paragraph-mark engine quirk:
# inserting/looking-up a truly zero-length Range right at a
# paragraph's end position (where the paragraph mark lives)
# confuses Bookmarks.Add. To avoid it, temporarily insert a unique
# marker string at that exact spot (which shifts the paragraph mark
# out of the way), anchor the bookmark at the now-safe position
# immediately in front of the marker, and finally delete the marker
# text again, leaving a clean zero-length "_GoBack" bookmark in the
# correct place.
$endPoint.InsertAfter("@@GOBACKMARK@@")

$bmAnchor = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmAnchor) | Out-Null

$markerRng = $d.Content
$markerRng.Find.Execute(
    "@@GOBACKMARK@@",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRng.Text = ""

Write-Host "Done"
